$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Babar%Kayani%xref no email%0,           Elliot%Onochie%xref no email%0,           Vijay%Patil%xref no email%0,           Fahima%Begum%xref no email%0,           Rory%Cuthbert%xref no email%0,           David%Ferguson%xref no email%0,           Jagmeet S.%Bhamra%xref no email%0,           Aadhar%Sharma%xref no email%0,           Peter%Bates%xref no email%0,           Fares S.%Haddad%xref no email%0]'
$ws.Range("E3").Value = '[Kenneth A.%Egol%NULL%0,                  Sanjit R.%Konda%NULL%0,                  Mackenzie L.%Bird%NULL%0,                  Nicket%Dedhia%NULL%0,                  Emma K.%Landes%NULL%0,                  Rachel A.%Ranson%NULL%0,                  Sara J.%Solasz%NULL%0,                  Vinay K.%Aggarwal%NULL%0,                  Joseph A.%Bosco%NULL%0,                  David L.%Furgiuele%NULL%0,                  Abhishek%Ganta%NULL%0,                  Jason%Gould%NULL%0,                  Thomas R.%Lyon%NULL%0,                  Toni M.%McLaurin%NULL%0,                  Nirmal C.%Tejwani%NULL%0,                  Joseph D.%Zuckerman%NULL%0,                  Philipp%Leucht%NULL%0]'
$ws.Range("E4").Value = '[Drake G.%LeBrun%NULL%0,                  Maxwell A.%Konnaris%NULL%0,                  Gregory C.%Ghahramani%NULL%0,                  Ajay%Premkumar%NULL%0,                  Chris J.%DeFrancesco%NULL%0,                  Jordan A.%Gruskay%NULL%0,                  Aleksey%Dvorzhinskiy%NULL%0,                  Milan S.%Sandhu%NULL%0,                  Elan M.%Goldwyn%NULL%0,                  Christopher L.%Mendias%NULL%0,                  William M.%Ricci%NULL%0]'
$ws.Range("E5").Value = '[Amit%Thakrar%NULL%1,                  Karen%Chui%NULL%1,                  Akhil%Kapoor%NULL%1,                  John%Hambidge%NULL%1]'
$ws.Range("E6").Value = '[Andrew J.%Hall%xref no email%0,           Nicholas D.%Clement%xref no email%0,           Luke%Farrow%xref no email%0,           Alasdair M. J.%MacLullich%xref no email%0,           Graham F.%Dall%xref no email%0,           Chloe E. H.%Scott%xref no email%0,           Paul J.%Jenkins%xref no email%0,           Timothy O.%White%xref no email%0,           Andrew D.%Duckworth%xref no email%0]'
$ws.Range("E7").Value = '[Karen%Chui%NULL%1,                  Amit%Thakrar%NULL%1,                  Shivakumar%Shankar%NULL%1]'
